$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$nl = [char]10

# ---------------------------------------------------------------------------
# Sheet1 ("scrape") header row updates
# ---------------------------------------------------------------------------
$audienceText = "Audience" + $nl + "(student, teacher, site administrator, district administrator, parent, counselor, implementation coordinator)"
$ws1.Range("F1").Value = $audienceText
$ws2.Range("I1").Value = $audienceText

$teacherAudience = "Teacher, Site Administrator, District Administrator, Counselor, Implementation Coordinator"

# ---------------------------------------------------------------------------
# Sheet1 data rows: re-point rows 3-7 to the "new" resources, normalize
# formatting on A4 / B6 / B7, then append row 8 for the new resource.
# ---------------------------------------------------------------------------

# Row 2 keeps its original resource (Evaluation Guide) - only the audience
# text (column F) needs the "site" wording.
$ws1.Range("F2").Value = $teacherAudience

# Row 3: Evaluation Design -> Evaluating Implementation
$ws1.Range("A3").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/Evaluating-Implementation"
$ws1.Range("B3").Value = "Evaluating Implementation"
$ws1.Range("F3").Value = $teacherAudience
$ws1.Range("A3").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A3"), $ws1.Range("A3").Value())

# Row 4: Outcome Measures -> Evaluation Design
$ws1.Range("A4").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/evaluation-design"
$ws1.Range("B4").Value = "Evaluation Design"
$ws1.Range("F4").Value = $teacherAudience
$ws1.Range("A4").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A4"), $ws1.Range("A4").Value())
# normalize A4 formatting to match the other rows in column A (drop its
# one-off style)
$ws1.Range("A3").Copy()
$ws1.Range("A4").PasteSpecial(-4122)

# Row 5: Using Data -> Using School Data
$ws1.Range("A5").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/using-school-data"
$ws1.Range("B5").Value = "Using School Data"
$ws1.Range("F5").Value = $teacherAudience
$ws1.Range("A5").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A5"), $ws1.Range("A5").Value())

# Row 6: Using Evaluation Findings -> Using Evaluation Results
$ws1.Range("A6").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/Evaluation-Guidehttp://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/using-evaluation-results"
$ws1.Range("B6").Value = "Using Evaluation Results"
$ws1.Range("F6").Value = $teacherAudience
$ws1.Range("A6").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A6"), $ws1.Range("A6").Value())
# normalize B6 formatting to match the other rows in column B
$ws1.Range("B2").Copy()
$ws1.Range("B6").PasteSpecial(-4122)

# Row 7: Evaluating Implementation -> Outcome Measures
$ws1.Range("A7").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/Evaluation-Guidehttp://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/Outcome-Measures"
$ws1.Range("B7").Value = "Outcome Measures"
$ws1.Range("F7").Value = $teacherAudience
$ws1.Range("A7").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A7"), $ws1.Range("A7").Value())
$ws1.Range("B2").Copy()
$ws1.Range("B7").PasteSpecial(-4122)

# Row 8 (new): More Evaluation Tools - copy row 7 formatting down first
$ws1.Range("A7:G7").Copy()
$ws1.Range("A8:G8").PasteSpecial(-4122)
$ws1.Range("A8").Value = "http://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/Evaluation-Guidehttp://www.secondstep.org/Kindergarten/Program-Coordinators/Second-Step-Kit/evaluation-guide/More-Evaluation-Tools"
$ws1.Range("B8").Value = "More Evaluation Tools"
$ws1.Range("C8").Value = $ws1.Range("C7").Value()
$ws1.Range("D8").Value = $ws1.Range("D7").Value()
$ws1.Range("E8").Value = $ws1.Range("E7").Value()
$ws1.Range("F8").Value = $teacherAudience
$ws1.Range("G8").Value = $ws1.Range("G7").Value()
$ws1.Hyperlinks.Add($ws1.Range("A8"), $ws1.Range("A8").Value())

# ---------------------------------------------------------------------------
# Row heights (sheet1)
# ---------------------------------------------------------------------------
$ws1.Rows.Item(1).RowHeight = 120.75
$ws1.Rows.Item(2).RowHeight = 86.25
$ws1.Rows.Item(3).RowHeight = 86.25
$ws1.Rows.Item(4).RowHeight = 86.25
$ws1.Rows.Item(5).RowHeight = 86.25
$ws1.Rows.Item(6).RowHeight = 120
$ws1.Rows.Item(7).RowHeight = 120
$ws1.Rows.Item(8).RowHeight = 120

# ---------------------------------------------------------------------------
# Sheet1 selection / view
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B8").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("File Properties") updates
# ---------------------------------------------------------------------------
$ws2.Range("I2").Value = $teacherAudience
$ws2.Range("H2").Value = $false
$ws2.Range("J2").Value = $true

$ws2.Rows.Item(1).RowHeight = 85.5
$ws2.Rows.Item(2).RowHeight = 69

$ws2.Activate()
$ws2.Range("I6").Select()

$ws1.Activate()
